$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two BTEC/Pearson logo pictures in the header story
# (the "default" header and the "first page" header both carry the same
# BTEC logo) and two in the footer story (default + first page, both
# carrying the Pearson logo). Each inline picture's display Name is being
# swapped to the *other* sibling file's name.

# Headers: both inline BTEC logo pictures get renamed to "image1.jpg"
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $rng = $hdr.Range
        for ($j = 1; $j -le $rng.InlineShapes.Count; $j++) {
            $shp = $rng.InlineShapes.Item($j)
            $shp.Name = "image1.jpg"
        }
    }
}

# Footers: both inline Pearson logo pictures get renamed to "image2.png"
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $rng = $ftr.Range
        for ($j = 1; $j -le $rng.InlineShapes.Count; $j++) {
            $shp = $rng.InlineShapes.Item($j)
            $shp.Name = "image2.png"
        }
    }
}

Write-Output "Renamed header/footer logo pictures."
